$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Cells.Item(28, 8).Value = 478.44446
$ws.Cells.Item(28, 9).Value = 299.77777
$ws.Cells.Item(28, 10).Value = 835.7778
$ws.Cells.Item(28, 11).Value = 299.77777
$ws.Cells.Item(28, 12).Value = 835.7778
$ws.Cells.Item(28, 13).Value = 185.22223
$ws.Cells.Item(28, 14).Value = -1805.7778
# Row 57
$ws.Cells.Item(57, 8).Value = 135181.86
$ws.Cells.Item(57, 10).Value = 135181.86
$ws.Cells.Item(57, 12).Value = 405545.58
$ws.Cells.Item(57, 14).Value = -406543.58
# Row 64
$ws.Cells.Item(64, 8).Value = 3720.4167
$ws.Cells.Item(64, 9).Value = 3080.75
$ws.Cells.Item(64, 10).Value = 4999.75
$ws.Cells.Item(64, 11).Value = 3080.75
$ws.Cells.Item(64, 12).Value = 4999.75
$ws.Cells.Item(64, 13).Value = -2832.75
$ws.Cells.Item(64, 14).Value = -5495.75
# Row 67
$ws.Cells.Item(67, 8).Value = 3720.4167
$ws.Cells.Item(67, 9).Value = 3080.75
$ws.Cells.Item(67, 10).Value = 4999.75
$ws.Cells.Item(67, 11).Value = 3080.75
$ws.Cells.Item(67, 12).Value = 4999.75
$ws.Cells.Item(67, 13).Value = -2222.75
$ws.Cells.Item(67, 14).Value = -6715.75
# Row 70
$ws.Cells.Item(70, 8).Value = 1361.5
$ws.Cells.Item(70, 9).Value = 2999
$ws.Cells.Item(70, 10).Value = 1244.5358
$ws.Cells.Item(70, 11).Value = 8997
$ws.Cells.Item(70, 12).Value = 3733.6074
$ws.Cells.Item(70, 13).Value = -8727
$ws.Cells.Item(70, 14).Value = -4273.607400000001
# Row 73
$ws.Cells.Item(73, 8).Value = 1361.5
$ws.Cells.Item(73, 9).Value = 2999
$ws.Cells.Item(73, 10).Value = 1244.5358
$ws.Cells.Item(73, 11).Value = 8997
$ws.Cells.Item(73, 12).Value = 3733.6074
$ws.Cells.Item(73, 13).Value = -8061
$ws.Cells.Item(73, 14).Value = -5605.607400000001
# Row 93
$ws.Cells.Item(93, 8).Value = 29062.666
$ws.Cells.Item(93, 10).Value = 29062.666
$ws.Cells.Item(93, 12).Value = 29062.666
$ws.Cells.Item(93, 14).Value = -34054.666
# Row 95
$ws.Cells.Item(95, 8).Value = 40311.5
$ws.Cells.Item(95, 10).Value = 40311.5
$ws.Cells.Item(95, 12).Value = 40311.5
$ws.Cells.Item(95, 14).Value = -45803.5
# Row 116
$ws.Cells.Item(116, 8).Value = 14677.6
$ws.Cells.Item(116, 10).Value = 20166.666
$ws.Cells.Item(116, 12).Value = 20166.666
$ws.Cells.Item(116, 14).Value = -27050.666
# Row 132
$ws.Cells.Item(132, 8).Value = 8458.229499999999
$ws.Cells.Item(132, 9).Value = 5777.706
$ws.Cells.Item(132, 11).Value = 17333.118
$ws.Cells.Item(132, 13).Value = -14803.118

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 1820.0454
$ws.Cells.Item(2, 9).Value = 1820.0454
$ws.Cells.Item(2, 11).Value = 1820.0454
$ws.Cells.Item(2, 13).Value = -1707.0454
# Row 63
$ws.Cells.Item(63, 8).Value = 1880.0714
$ws.Cells.Item(63, 10).Value = 2430
$ws.Cells.Item(63, 12).Value = 2430
$ws.Cells.Item(63, 14).Value = -3802
# Row 66
$ws.Cells.Item(66, 8).Value = 1880.0714
$ws.Cells.Item(66, 10).Value = 2430
$ws.Cells.Item(66, 12).Value = 12150
$ws.Cells.Item(66, 14).Value = -19014
# Row 102
$ws.Cells.Item(102, 8).Value = 200005490
$ws.Cells.Item(102, 9).Value = 6867.5
$ws.Cells.Item(102, 10).Value = 1000000000
$ws.Cells.Item(102, 11).Value = 6867.5
$ws.Cells.Item(102, 12).Value = 1000000000
$ws.Cells.Item(102, 13).Value = -5245.5
$ws.Cells.Item(102, 14).Value = -1000003244
# Row 116
$ws.Cells.Item(116, 8).Value = 1820.0454
$ws.Cells.Item(116, 9).Value = 1820.0454
$ws.Cells.Item(116, 11).Value = 1820.0454
$ws.Cells.Item(116, 13).Value = 473.9546
# Row 132
$ws.Cells.Item(132, 8).Value = 2012.0488
$ws.Cells.Item(132, 9).Value = 2012.0488
$ws.Cells.Item(132, 11).Value = 6036.1464
$ws.Cells.Item(132, 13).Value = -3506.1464

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 1820.0454
$ws.Cells.Item(3, 9).Value = 1820.0454
$ws.Cells.Item(3, 11).Value = 1820.0454
$ws.Cells.Item(3, 13).Value = -1706.0454
# Row 60
$ws.Cells.Item(60, 8).Value = 40000
$ws.Cells.Item(60, 9).Value = 40000
$ws.Cells.Item(60, 10).Value = 0
$ws.Cells.Item(60, 11).Value = 40000
$ws.Cells.Item(60, 12).Value = 0
$ws.Cells.Item(60, 13).Value = -39401
$ws.Cells.Item(60, 14).Value = $null

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 74
$ws.Cells.Item(74, 8).Value = 68380.836
$ws.Cells.Item(74, 9).Value = 10285
$ws.Cells.Item(74, 10).Value = 80000
$ws.Cells.Item(74, 11).Value = 10285
$ws.Cells.Item(74, 12).Value = 80000
$ws.Cells.Item(74, 13).Value = -9411
$ws.Cells.Item(74, 14).Value = -81748
# Row 77
$ws.Cells.Item(77, 8).Value = 68380.836
$ws.Cells.Item(77, 9).Value = 10285
$ws.Cells.Item(77, 10).Value = 80000
$ws.Cells.Item(77, 11).Value = 30855
$ws.Cells.Item(77, 12).Value = 240000
$ws.Cells.Item(77, 13).Value = -26487
$ws.Cells.Item(77, 14).Value = -248736
# Row 107
$ws.Cells.Item(107, 8).Value = 775.4318
$ws.Cells.Item(107, 9).Value = 634.1111
$ws.Cells.Item(107, 11).Value = 634.1111
$ws.Cells.Item(107, 13).Value = 1285.8889

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Cells.Item(3, 8).Value = 2932.375
$ws.Cells.Item(3, 9).Value = 2932.375
$ws.Cells.Item(3, 11).Value = 8797.125
$ws.Cells.Item(3, 13).Value = -8685.125
# Row 5
$ws.Cells.Item(5, 8).Value = 747.125
$ws.Cells.Item(5, 10).Value = 245
$ws.Cells.Item(5, 12).Value = 735
$ws.Cells.Item(5, 14).Value = -959
# Row 118
$ws.Cells.Item(118, 8).Value = 600
$ws.Cells.Item(118, 9).Value = 600
$ws.Cells.Item(118, 11).Value = 1800
$ws.Cells.Item(118, 13).Value = -557
# Row 135
$ws.Cells.Item(135, 8).Value = 747.125
$ws.Cells.Item(135, 10).Value = 245
$ws.Cells.Item(135, 12).Value = 2205
$ws.Cells.Item(135, 14).Value = -7275

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 41
$ws.Cells.Item(41, 8).Value = 8365.5
$ws.Cells.Item(41, 9).Value = 8438.6
$ws.Cells.Item(41, 11).Value = 8438.6
$ws.Cells.Item(41, 13).Value = -8083.6
# Row 57
$ws.Cells.Item(57, 8).Value = 3027.5
$ws.Cells.Item(57, 9).Value = 3027.5
$ws.Cells.Item(57, 10).Value = 0
$ws.Cells.Item(57, 11).Value = 3027.5
$ws.Cells.Item(57, 12).Value = 0
$ws.Cells.Item(57, 13).Value = -2207.5
$ws.Cells.Item(57, 14).Value = $null
# Row 105
$ws.Cells.Item(105, 8).Value = 42467.4
$ws.Cells.Item(105, 10).Value = 42467.4
$ws.Cells.Item(105, 12).Value = 42467.4
$ws.Cells.Item(105, 14).Value = -49455.4
# Row 132
$ws.Cells.Item(132, 8).Value = 4440.206
$ws.Cells.Item(132, 9).Value = 4131.593
$ws.Cells.Item(132, 10).Value = 5630.5713
$ws.Cells.Item(132, 11).Value = 12394.779
$ws.Cells.Item(132, 12).Value = 16891.7139
$ws.Cells.Item(132, 13).Value = -9864.778999999999
$ws.Cells.Item(132, 14).Value = -21951.7139

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells.Item(22, 8).Value = 994.8
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 14).Value = $null
# Row 27
$ws.Cells.Item(27, 8).Value = 994.8
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 12).Value = 0
$ws.Cells.Item(27, 14).Value = $null
# Row 115
$ws.Cells.Item(115, 8).Value = 150000
$ws.Cells.Item(115, 10).Value = 150000
$ws.Cells.Item(115, 12).Value = 150000
$ws.Cells.Item(115, 14).Value = -152350
# Row 132
$ws.Cells.Item(132, 8).Value = 36820.35
$ws.Cells.Item(132, 9).Value = 36820.35
$ws.Cells.Item(132, 11).Value = 110461.05
$ws.Cells.Item(132, 13).Value = -107931.05

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Cells.Item(132, 8).Value = 2239.8235
$ws.Cells.Item(132, 9).Value = 2269.24
$ws.Cells.Item(132, 11).Value = 6807.719999999999
$ws.Cells.Item(132, 13).Value = -4277.719999999999
